$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F33").Value = '18_hazards_to_humans_and_domestic_animals'
$ws.Range("F34").Value = '18_hazards_to_humans_and_domestic_animals'
$ws.Range("F35").Value = '18_hazards_to_humans_and_domestic_animals'
$ws.Range("F36").Value = 'application instructions || env warning - species'
$ws.Range("F55").Value = '154_pesticide_storage'
$ws.Range("F103").Value = '18_hazards_to_humans_and_domestic_animals'
$ws.Range("F104").Value = 'application instructions || env warning - species'
$ws.Range("F105").Value = 'ppe'
$ws.Range("F106").Value = 'ppe'
$ws.Range("F107").Value = 'ppe'
$ws.Range("F108").Value = 'ppe'
$ws.Range("F112").Value = 'mixing'
$ws.Range("F115").Value = 'use restrictions'
$ws.Range("F116").Value = 'use restrictions'
$ws.Range("F119").Value = 'mixing'
$ws.Range("F130").Value = 'application instructions'
$ws.Range("F131").Value = 'use restrictions'
$ws.Range("F132").Value = 'application instructions'
$ws.Range("F133").Value = 'application instructions'
$ws.Range("F134").Value = 'application instructions'
$ws.Range("F135").Value = 'application instructions'
$ws.Range("F136").Value = 'application instructions'
$ws.Range("F139").Value = 'application instructions'
$ws.Range("F142").Value = 'application instructions'
$ws.Range("F143").Value = 'application instructions'
$ws.Range("F144").Value = 'application instructions'
$ws.Range("F145").Value = 'application instructions'
$ws.Range("F146").Value = 'application instructions'
$ws.Range("F150").Value = 'application instructions'
$ws.Range("F154").Value = 'application instructions'
$ws.Range("F158").Value = 'use restrictions'
$ws.Range("F159").Value = 'application instructions'
$ws.Range("F163").Value = 'application instructions'
$ws.Range("F164").Value = 'application instructions'
$ws.Range("F170").Value = 'application instructions'
$ws.Range("F171").Value = 'application instructions'
$ws.Range("F172").Value = 'application instructions'
$ws.Range("F176").Value = 'application instructions'
$ws.Range("F178").Value = 'application instructions'
$ws.Range("F179").Value = 'application instructions'
$ws.Range("F188").Value = 'use restrictions || application instructions'
$ws.Range("F192").Value = 'use restrictions'
$ws.Range("F193").Value = 'use restrictions'
$ws.Range("F194").Value = 'use restrictions'
$ws.Range("F196").Value = 'application instructions'
$ws.Range("F197").Value = 'application instructions'
$ws.Range("F198").Value = 'application instructions'
$ws.Range("F200").Value = '154_pesticide_storage'
